$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("DataFormat")

# --- Sheet "Data": add Title / Department / ActivityListNames headers ---
$ws1.Cells.Item(1, 15).Value = "Title"
$ws1.Cells.Item(1, 16).Value = "Department"
$ws1.Cells.Item(1, 17).Value = "ActivityListNames"

# --- Sheet "DataFormat": duplicate header row (row 3) ---
$ws2.Cells.Item(3, 15).Value = "Title"
$ws2.Cells.Item(3, 16).Value = "Department"
$ws2.Cells.Item(3, 17).Value = "ActivityListNames"

# --- Sheet "DataFormat": sample data row (row 4) ---
$ws2.Cells.Item(4, 15).Value = "Manager"
$ws2.Cells.Item(4, 16).Value = "Accounts Payable"
$ws2.Cells.Item(4, 17).Value = "Welcome Dinner, Universal Studios, Breakout Session 1701"

# --- Sheet "DataFormat": field-type legend rows 21-23 ---
$ws2.Cells.Item(21, 1).Value = "Title"
$ws2.Cells.Item(21, 1).Font.Bold = $true
$ws2.Cells.Item(21, 2).Value = "String"

$ws2.Cells.Item(22, 1).Value = "Department"
$ws2.Cells.Item(22, 1).Font.Bold = $true
$ws2.Cells.Item(22, 2).Value = "String"

$ws2.Cells.Item(23, 1).Value = "ActivityListNames"
$ws2.Cells.Item(23, 1).Font.Bold = $true
$ws2.Cells.Item(23, 2).Value = 'String, separated by commas (",") for future parsing'

# --- Selections to match the final saved state ---
[void]$ws1.Range("Q2").Select()

[void]$ws2.Activate()
[void]$ws2.Range("B24").Select()
